# "add NPC move type"
#
# - W1 header "DropPack" is renamed to "DropPackList" and shifts out to Y1.
# - New headers "MoveType" (W1) and "AtkDis" (X1) are inserted.
#   The "AtkDis" header cell gets a dedicated font (family 3) and the run
#   "tkDis" (characters 2-6) carries that font explicitly as rich text,
#   while the leading "A" keeps the default font.
# - The new MoveType / AtkDis columns are populated for the 5 data rows.
# - Columns V:Y are resized.
# - The sheet view scrolls right and re-selects X10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -------------------------------------------------------------
# Write Y1 first so the shared-string slot the original "DropPack" header
# occupied gets renamed in place to "DropPackList" rather than orphaned.
$ws.Range("Y1").Value = "DropPackList"
$ws.Range("W1").Value = "MoveType"
$ws.Range("X1").Value = "AtkDis"

# Give the whole X1 cell the new (family 3) font first...
$ws.Range("X1").Font.Name = "宋体"
$ws.Range("X1").Font.Size = 11
$ws.Range("X1").Font.ColorIndex = 1
$ws.Range("X1").Font.Family = 3

# ...then re-apply it explicitly to the "tkDis" run so it is recorded as
# its own rich-text run (the leading "A" is left at the default font).
$atkDisRun = $ws.Range("X1").Characters(2, 5)
$atkDisRun.Font.Name = "宋体"
$atkDisRun.Font.Size = 11
$atkDisRun.Font.ColorIndex = 1
$atkDisRun.Font.Family = 3

# --- Data ------------------------------------------------------------------
$moveType = @(2, 2, 2, 0, 2)
$atkDis = @(20, 20, 20, 20, 20)
for ($i = 0; $i -lt 5; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 23).Value = $moveType[$i]
    $ws.Cells.Item($row, 24).Value = $atkDis[$i]
}

# --- Column widths -----------------------------------------------------
$ws.Range("V1").ColumnWidth = 25
$ws.Range("W1:X1").ColumnWidth = 25
$ws.Range("Y1").ColumnWidth = 13.875

# --- View state ----------------------------------------------------------
$ws.Range("X10").Select()
$excel.ActiveWindow.ScrollColumn = 11
